# Append new "Success - <timestamp>" shared-string values to the test-run
# status columns on each sheet of the Jostens QA test-data workbook.

$wb = $excel.ActiveWorkbook

# Login sheet
$ws = $wb.Worksheets.Item("Login")
$ws.Range("G2").Value = "Success - 2020/12/19 16:58:57"
$ws.Range("G3").Value = "Success - 2020/12/19 16:59:00"

# School Search sheet
$ws = $wb.Worksheets.Item("School Search")
$ws.Range("C2").Value = "Success - 2020/12/19 16:59:03"
$ws.Range("C3").Value = "Success - 2020/12/19 16:59:06"

# Product Search sheet
$ws = $wb.Worksheets.Item("Product Search")
$ws.Range("K2").Value = "Success - 2020/12/19 16:59:32"
$ws.Range("K3").Value = "Success - 2020/12/19 16:59:54"
$ws.Range("K4").Value = "Success - 2020/12/19 17:00:15"

# Shopping Cart sheet
$ws = $wb.Worksheets.Item("Shopping Cart")
$ws.Range("G2").Value = "Success - 2020/12/19 17:00:17"
$ws.Range("G3").Value = "Success - 2020/12/19 17:00:17"
$ws.Range("G4").Value = "Success - 2020/12/19 17:00:17"

# Checkout sheet
$ws = $wb.Worksheets.Item("Checkout")
$ws.Range("P2").Value = "Success - 2020/12/19 17:00:27"
$ws.Range("P3").Value = "Success - 2020/12/19 17:00:38"
$ws.Range("P4").Value = "Success - 2020/12/19 17:00:46"

# Payment sheet (C2 carries a quote-prefix / text style - preserve it by
# typing the new value with a leading apostrophe, same as the original author)
$ws = $wb.Worksheets.Item("Payment")
$ws.Range("C2").Value = "'Success - 2020/12/19 17:00:56"
